# HIVE TEAMS.docx -- apply Russian translations to the remaining
# English paragraphs (the rest of the document was already translated).
#
# Paragraphs are addressed by their 1-based Word COM index
# ($d.Paragraphs.Item(N)) and the whole-paragraph text is replaced via
# Range.Text assignment (this performs a plain textual substitution that
# keeps the existing run's formatting and avoids Word's smart-quote /
# smart-dash AutoCorrect that Find.Execute's replacement text can trigger).

$d = $word.ActiveDocument

function Set-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Text = $newText
}

Set-ParaText 55 "Команда Hive: Разработка"
Set-ParaText 56 "Разработка экосистемы SmartCash, сервисы и приложения."
Set-ParaText 58 "Координатор Hive"
Set-ParaText 59 "Создатель Dash N Drink Soda Machine и SmartCash POS."
Set-ParaText 61 "Разработчик"
Set-ParaText 65 "Разработчик"
Set-ParaText 67 "Разработчик"
Set-ParaText 68 "Команда Hive: Продвижение (II)"
Set-ParaText 69 "Работа над расширением сообщества в Латинской Америке."
Set-ParaText 71 "Координатор Hive"
Set-ParaText 73 "Продвижение"
Set-ParaText 75 "Продвижение"
Set-ParaText 77 "Продвижение"

# Paragraph 78 has three runs: " HIVE TEAM: SUPPORT " + "&" + " WEB".
# Only the first and third runs change; the "&" run is left untouched,
# so target each piece individually via a range-scoped Find.Execute.
$p78 = $d.Paragraphs.Item(78)
$r78 = $p78.Range
$r78.Find.Execute(" HIVE TEAM: SUPPORT ", $true, $false, $false, $false, $false, $true, 1, $false, "Команда Hive: Поддержка ", 2)
$p78b = $d.Paragraphs.Item(78)
$r78b = $p78b.Range
$r78b.Find.Execute(" WEB", $true, $false, $false, $false, $false, $true, 1, $false, "Web", 2)

Set-ParaText 79 "Интеграции и поддержка пользователей."
Set-ParaText 81 "Координатор Hive"
Set-ParaText 82 'Alex – настоящий "человек всех профессий". Его конёк - технологии, графика, работа с инфраструктурой сайтов и веб-дизайн.'
Set-ParaText 84 "Финансовый управляющий"
Set-ParaText 86 "Поддержка"
Set-ParaText 88 "Поддержка"
Set-ParaText 90 "Поддержка"
Set-ParaText 92 "Вице-координатор"
Set-ParaText 94 "Видео-гуру"
Set-ParaText 96 "Юридические вопросы"
Set-ParaText 99 "Хотите присоединиться?"
Set-ParaText 100 "SmartHive – это место, где ваши таланты и способности будут востребованы. Присоединяйтесь!  "
Set-ParaText 101 "Мы считаем, что постоянные команды — это то, что приводит к коррупции и неэффективности. Мы хотим идти по собственному пути и поэтому создали децентрализованную организационную модель, основанную на принципах жизни муравьёв и пчелиных колоний."
Set-ParaText 102 "Создание и поддержание такой структуры управления требует особого подхода, поэтому мы разработали две концепции – SmartHive и Hive Structuring Teams (HST).  SmartHive дает возможность любому держателю монет голосовать за проекты и идеи, представленные сообществом.  Благодаря SmartHive каждый участник способен проявлять себя – выдвигать свои идеи на голосование, участвовать в обсуждениях, способствовать росту сообщества, а также голосовать за другие предложения."

Write-Output "Done"
